$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 21 de Octubre de 2020 a las 21:35"

# Update per-country COVID statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes). Curazao overtakes San Marino in
# total cases, so their two rows (170-171) swap places as well as getting new figures.
$ws.Cells.Item(4,2).Value2 = 8549673
$ws.Cells.Item(4,3).Value2 = 28723
$ws.Cells.Item(4,4).Value2 = 5557126
$ws.Cells.Item(4,5).Value2 = 2765828
$ws.Cells.Item(4,7).Value2 = 535
$ws.Cells.Item(4,8).Value2 = 226719
$ws.Cells.Item(5,2).Value2 = 7705158
$ws.Cells.Item(5,3).Value2 = 56000
$ws.Cells.Item(5,4).Value2 = 6871898
$ws.Cells.Item(5,5).Value2 = 716607
$ws.Cells.Item(5,7).Value2 = 703
$ws.Cells.Item(5,8).Value2 = 116653
$ws.Cells.Item(11,4).Value2 = 107652
$ws.Cells.Item(11,5).Value2 = 815721
$ws.Cells.Item(21,2).Value2 = 389561
$ws.Cells.Item(21,3).Value2 = 8663
$ws.Cells.Item(21,5).Value2 = 77473
$ws.Cells.Item(21,7).Value2 = 33
$ws.Cells.Item(21,8).Value2 = 9988
$ws.Cells.Item(28,2).Value2 = 307259
$ws.Cells.Item(28,3).Value2 = 1097
$ws.Cells.Item(28,4).Value2 = 285411
$ws.Cells.Item(28,5).Value2 = 19557
$ws.Cells.Item(28,7).Value2 = 13
$ws.Cells.Item(28,8).Value2 = 2291
$ws.Cells.Item(35,5).Value2 = 27867
$ws.Cells.Item(35,7).Value2 = 52
$ws.Cells.Item(35,8).Value2 = 3079
$ws.Cells.Item(62,2).Value2 = 65577
$ws.Cells.Item(62,3).Value2 = 1241
$ws.Cells.Item(62,4).Value2 = 30470
$ws.Cells.Item(62,5).Value2 = 34571
$ws.Cells.Item(62,7).Value2 = 5
$ws.Cells.Item(62,8).Value2 = 536
$ws.Cells.Item(64,4).Value2 = 56798
$ws.Cells.Item(64,5).Value2 = 3707
$ws.Cells.Item(72,2).Value2 = 48628
$ws.Cells.Item(72,3).Value2 = 499
$ws.Cells.Item(72,4).Value2 = 41935
$ws.Cells.Item(72,5).Value2 = 6266
$ws.Cells.Item(72,7).Value2 = 6
$ws.Cells.Item(72,8).Value2 = 427
$ws.Cells.Item(76,2).Value2 = 45892
$ws.Cells.Item(76,3).Value2 = 1442
$ws.Cells.Item(76,5).Value2 = 40120
$ws.Cells.Item(76,7).Value2 = 29
$ws.Cells.Item(76,8).Value2 = 740
$ws.Cells.Item(104,2).Value2 = 12406
$ws.Cells.Item(104,3).Value2 = 39
$ws.Cells.Item(104,4).Value2 = 10587
$ws.Cells.Item(104,5).Value2 = 1686
$ws.Cells.Item(104,7).Value2 = 1
$ws.Cells.Item(104,8).Value2 = 133
$ws.Cells.Item(113,2).Value2 = 8979
$ws.Cells.Item(113,3).Value2 = 3
$ws.Cells.Item(113,5).Value2 = 1445
$ws.Cells.Item(119,2).Value2 = 8033
$ws.Cells.Item(119,3).Value2 = 132
$ws.Cells.Item(119,4).Value2 = 6835
$ws.Cells.Item(119,5).Value2 = 1108
$ws.Cells.Item(119,7).Value2 = 3
$ws.Cells.Item(119,8).Value2 = 90
$ws.Cells.Item(125,2).Value2 = 5864
$ws.Cells.Item(125,3).Value2 = 3
$ws.Cells.Item(125,4).Value2 = 4762
$ws.Cells.Item(125,5).Value2 = 919
$ws.Cells.Item(125,7).Value2 = 1
$ws.Cells.Item(125,8).Value2 = 183
$ws.Cells.Item(126,2).Value2 = 5805
$ws.Cells.Item(126,3).Value2 = 5
$ws.Cells.Item(126,4).Value2 = 5453
$ws.Cells.Item(126,5).Value2 = 236
$ws.Cells.Item(128,2).Value2 = 5512
$ws.Cells.Item(128,3).Value2 = 13
$ws.Cells.Item(128,4).Value2 = 5387
$ws.Cells.Item(128,5).Value2 = 64
$ws.Cells.Item(138,2).Value2 = 5012
$ws.Cells.Item(138,3).Value2 = 16
$ws.Cells.Item(138,4).Value2 = 4798
$ws.Cells.Item(138,5).Value2 = 180
$ws.Cells.Item(151,2).Value2 = 3428
$ws.Cells.Item(151,3).Value2 = 17
$ws.Cells.Item(151,4).Value2 = 2599
$ws.Cells.Item(151,5).Value2 = 697
$ws.Cells.Item(165,2).Value2 = 1404
$ws.Cells.Item(165,3).Value2 = 5
$ws.Cells.Item(165,4).Value2 = 1221
$ws.Cells.Item(165,5).Value2 = 87
$ws.Cells.Item(165,7).Value2 = 3
$ws.Cells.Item(165,8).Value2 = 96
$ws.Cells.Item(166,2).Value2 = 1385
$ws.Cells.Item(166,3).Value2 = 1
$ws.Cells.Item(166,4).Value2 = 1276
$ws.Cells.Item(166,5).Value2 = 27
$ws.Cells.Item(170,1).Value2 = "Curazao"
$ws.Cells.Item(170,2).Value2 = 785
$ws.Cells.Item(170,3).Value2 = 26
$ws.Cells.Item(170,4).Value2 = 479
$ws.Cells.Item(170,5).Value2 = 305
$ws.Cells.Item(170,8).Value2 = 1
$ws.Cells.Item(171,1).Value2 = "San Marino"
$ws.Cells.Item(171,2).Value2 = 774
$ws.Cells.Item(171,3).Value2 = 8
$ws.Cells.Item(171,4).Value2 = 690
$ws.Cells.Item(171,5).Value2 = 42
$ws.Cells.Item(171,8).Value2 = 42
$ws.Cells.Item(189,2).Value2 = 273
$ws.Cells.Item(189,3).Value2 = 2
$ws.Cells.Item(189,4).Value2 = 231
$ws.Cells.Item(189,5).Value2 = 40